$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1900
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 1625
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 1625
$ws.Range("M32").Value = -2674
$ws.Range("N32").Value = -2277
$ws.Range("H40").Value = 106705.2
$ws.Range("J40").Value = 3875
$ws.Range("L40").Value = 3875
$ws.Range("N40").Value = -4225
$ws.Range("H86").Value = 281251260
$ws.Range("I86").Value = 250001460
$ws.Range("K86").Value = 250001460
$ws.Range("M86").Value = -250000337
$ws.Range("H89").Value = 281251260
$ws.Range("I89").Value = 250001460
$ws.Range("K89").Value = 1250007300
$ws.Range("M89").Value = -1250001684
$ws.Range("H96").Value = 1836.3125
$ws.Range("I96").Value = 1409.5385
$ws.Range("J96").Value = 3685.6667
$ws.Range("K96").Value = 4228.6155
$ws.Range("L96").Value = 11057.0001
$ws.Range("M96").Value = -2855.6155
$ws.Range("N96").Value = -13803.0001
$ws.Range("H106").Value = 4764720
$ws.Range("I106").Value = 4764720
$ws.Range("K106").Value = 4764720
$ws.Range("M106").Value = -4764089
$ws.Range("H112").Value = 1238184.4
$ws.Range("J112").Value = 1758508.8
$ws.Range("L112").Value = 5275526.4
$ws.Range("N112").Value = -5277742.4
$ws.Range("H125").Value = 2181.3076
$ws.Range("I125").Value = 998.4
$ws.Range("K125").Value = 8985.6
$ws.Range("M125").Value = -6525.6
$ws.Range("H135").Value = 1062.9445
$ws.Range("I135").Value = 1092.6875
$ws.Range("K135").Value = 9834.1875
$ws.Range("M135").Value = -7299.1875
$ws.Range("H137").Value = 2506.158
$ws.Range("I137").Value = 2588.3333
$ws.Range("J137").Value = 2365.2856
$ws.Range("K137").Value = 7764.999899999999
$ws.Range("L137").Value = 7095.8568
$ws.Range("M137").Value = -5214.999899999999
$ws.Range("N137").Value = -12195.8568
$ws.Range("H138").Value = 2547.37
$ws.Range("J138").Value = 3480.9333
$ws.Range("L138").Value = 10442.7999
$ws.Range("N138").Value = -20722.7999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4681.1665
$ws.Range("I2").Value = 4727.5
$ws.Range("K2").Value = 4727.5
$ws.Range("M2").Value = -4614.5
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 20524
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 2286.5715
$ws.Range("I61").Value = 2336.2856
$ws.Range("K61").Value = 2336.2856
$ws.Range("M61").Value = -2124.2856
$ws.Range("H111").Value = 82024.11
$ws.Range("J111").Value = 82277.25
$ws.Range("L111").Value = 82277.25
$ws.Range("N111").Value = -90457.25
$ws.Range("H116").Value = 4681.1665
$ws.Range("I116").Value = 4727.5
$ws.Range("K116").Value = 4727.5
$ws.Range("M116").Value = -2433.5
$ws.Range("H121").Value = 58850.6
$ws.Range("J121").Value = 58850.6
$ws.Range("L121").Value = 58850.6
$ws.Range("N121").Value = -62344.6
$ws.Range("H122").Value = 3694.2654
$ws.Range("I122").Value = 2691.9688
$ws.Range("K122").Value = 8075.9064
$ws.Range("M122").Value = -5625.9064
$ws.Range("H136").Value = 2286.5715
$ws.Range("I136").Value = 2336.2856
$ws.Range("K136").Value = 7008.8568
$ws.Range("M136").Value = -4458.8568
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4681.1665
$ws.Range("I3").Value = 4727.5
$ws.Range("K3").Value = 4727.5
$ws.Range("M3").Value = -4613.5
$ws.Range("H5").Value = 339.9
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 899.5
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 899.5
$ws.Range("M5").Value = -87
$ws.Range("N5").Value = -1125.5
$ws.Range("H20").Value = 38350.645
$ws.Range("I20").Value = 57579.332
$ws.Range("K20").Value = 57579.332
$ws.Range("M20").Value = -57332.332
$ws.Range("H62").Value = 128000
$ws.Range("J62").Value = 128000
$ws.Range("L62").Value = 128000
$ws.Range("N62").Value = -129372
$ws.Range("H65").Value = 128000
$ws.Range("J65").Value = 128000
$ws.Range("L65").Value = 384000
$ws.Range("N65").Value = -390864
$ws.Range("H94").Value = 913.6667
$ws.Range("I94").Value = 829.71875
$ws.Range("K94").Value = 829.71875
$ws.Range("M94").Value = -378.71875
$ws.Range("H99").Value = 2554.5881
$ws.Range("I99").Value = 1935.2
$ws.Range("K99").Value = 1935.2
$ws.Range("M99").Value = -437.2
$ws.Range("H134").Value = 1663093.4
$ws.Range("I134").Value = 2166061.5
$ws.Range("K134").Value = 6498184.5
$ws.Range("M134").Value = -6495649.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3928.275
$ws.Range("I31").Value = 1755.1538
$ws.Range("K31").Value = 1755.1538
$ws.Range("M31").Value = -1460.1538
$ws.Range("H34").Value = 3928.275
$ws.Range("I34").Value = 1755.1538
$ws.Range("K34").Value = 1755.1538
$ws.Range("M34").Value = -1553.1538
$ws.Range("H58").Value = 2325.889
$ws.Range("I58").Value = 2227.4119
$ws.Range("K58").Value = 2227.4119
$ws.Range("M58").Value = -2024.4119
$ws.Range("H62").Value = 4624.75
$ws.Range("I62").Value = 4750
$ws.Range("J62").Value = 4499.5
$ws.Range("K62").Value = 4750
$ws.Range("L62").Value = 4499.5
$ws.Range("M62").Value = -4126
$ws.Range("N62").Value = -5747.5
$ws.Range("H65").Value = 4624.75
$ws.Range("I65").Value = 4750
$ws.Range("J65").Value = 4499.5
$ws.Range("K65").Value = 23750
$ws.Range("L65").Value = 22497.5
$ws.Range("M65").Value = -20630
$ws.Range("N65").Value = -28737.5
$ws.Range("H98").Value = 60929.5
$ws.Range("J98").Value = 60929.5
$ws.Range("L98").Value = 60929.5
$ws.Range("N98").Value = -65421.5
$ws.Range("H99").Value = 1748.75
$ws.Range("I99").Value = 1748.75
$ws.Range("K99").Value = 1748.75
$ws.Range("M99").Value = -250.75
$ws.Range("H105").Value = 1860.5
$ws.Range("I105").Value = 1141.6666
$ws.Range("J105").Value = 2938.75
$ws.Range("K105").Value = 1141.6666
$ws.Range("L105").Value = 2938.75
$ws.Range("M105").Value = 605.3334
$ws.Range("N105").Value = -6432.75
$ws.Range("H110").Value = 79975.664
$ws.Range("J110").Value = 79975.664
$ws.Range("L110").Value = 79975.664
$ws.Range("N110").Value = -88155.664
$ws.Range("H126").Value = 1748.75
$ws.Range("I126").Value = 1748.75
$ws.Range("K126").Value = 5246.25
$ws.Range("M126").Value = -2776.25
$ws.Range("H134").Value = 39268.445
$ws.Range("I134").Value = 49536.715
$ws.Range("K134").Value = 148610.145
$ws.Range("M134").Value = -146075.145
$ws.Range("H136").Value = 2325.889
$ws.Range("I136").Value = 2227.4119
$ws.Range("K136").Value = 6682.2357
$ws.Range("M136").Value = -4132.2357
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 7950.467
$ws.Range("I14").Value = 7950.467
$ws.Range("K14").Value = 23851.401
$ws.Range("M14").Value = -23678.401
$ws.Range("H20").Value = 2501
$ws.Range("I20").Value = 2501
$ws.Range("K20").Value = 7503
$ws.Range("M20").Value = -7276
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H80").Value = 2464
$ws.Range("I80").Value = 2336.6667
$ws.Range("J80").Value = 2591.3333
$ws.Range("K80").Value = 2336.6667
$ws.Range("L80").Value = 2591.3333
$ws.Range("M80").Value = -1338.6667
$ws.Range("N80").Value = -4587.3333
$ws.Range("H83").Value = 2464
$ws.Range("I83").Value = 2336.6667
$ws.Range("J83").Value = 2591.3333
$ws.Range("K83").Value = 11683.3335
$ws.Range("L83").Value = 12956.6665
$ws.Range("M83").Value = -6691.333500000001
$ws.Range("N83").Value = -22940.6665
$ws.Range("H97").Value = 1081.5862
$ws.Range("I97").Value = 885.56525
$ws.Range("J97").Value = 1833
$ws.Range("K97").Value = 885.56525
$ws.Range("L97").Value = 1833
$ws.Range("M97").Value = -389.56525
$ws.Range("N97").Value = -2825
$ws.Range("H112").Value = 110999
$ws.Range("J112").Value = 110999
$ws.Range("L112").Value = 110999
$ws.Range("N112").Value = -113215
$ws.Range("H122").Value = 4536.9473
$ws.Range("I122").Value = 3997
$ws.Range("J122").Value = 6048.8
$ws.Range("K122").Value = 11991
$ws.Range("L122").Value = 18146.4
$ws.Range("M122").Value = -9541
$ws.Range("N122").Value = -23046.4
$ws.Range("H126").Value = 3515.5
$ws.Range("I126").Value = 3431
$ws.Range("K126").Value = 10293
$ws.Range("M126").Value = -7823
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8160.1304
$ws.Range("I7").Value = 7360.3076
$ws.Range("K7").Value = 7360.3076
$ws.Range("M7").Value = -7248.3076
$ws.Range("H22").Value = 1981.76
$ws.Range("I22").Value = 1818.3334
$ws.Range("J22").Value = 2226.9
$ws.Range("K22").Value = 1818.3334
$ws.Range("L22").Value = 2226.9
$ws.Range("M22").Value = -1523.3334
$ws.Range("N22").Value = -2816.9
$ws.Range("H27").Value = 1981.76
$ws.Range("I27").Value = 1818.3334
$ws.Range("J27").Value = 2226.9
$ws.Range("K27").Value = 1818.3334
$ws.Range("L27").Value = 2226.9
$ws.Range("M27").Value = -1711.3334
$ws.Range("N27").Value = -2440.9
$ws.Range("H122").Value = 9957.808000000001
$ws.Range("I122").Value = 10026.789
$ws.Range("J122").Value = 9770.571
$ws.Range("K122").Value = 30080.367
$ws.Range("L122").Value = 29311.713
$ws.Range("M122").Value = -27630.367
$ws.Range("N122").Value = -34211.713
$ws.Range("H126").Value = 8160.1304
$ws.Range("I126").Value = 7360.3076
$ws.Range("K126").Value = 22080.9228
$ws.Range("M126").Value = -19610.9228
$ws.Range("H132").Value = 56473.35
$ws.Range("I132").Value = 68716.75
$ws.Range("K132").Value = 206150.25
$ws.Range("M132").Value = -203620.25
$ws.Range("H136").Value = 2378.24
$ws.Range("I136").Value = 1810.5883
$ws.Range("K136").Value = 5431.7649
$ws.Range("M136").Value = -2881.7649
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 30056
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H100").Value = 1022.5
$ws.Range("I100").Value = 1136.6875
$ws.Range("J100").Value = 794.125
$ws.Range("K100").Value = 2273.375
$ws.Range("L100").Value = 1588.25
$ws.Range("M100").Value = -1732.375
$ws.Range("N100").Value = -2670.25
$ws.Range("H132").Value = 2556.6064
$ws.Range("I132").Value = 2255.1555
$ws.Range("K132").Value = 6765.4665
$ws.Range("M132").Value = -4235.4665
$ws.Range("H136").Value = 42741.8
$ws.Range("I136").Value = 2572.4546
$ws.Range("J136").Value = 74303.42999999999
$ws.Range("K136").Value = 7717.3638
$ws.Range("L136").Value = 222910.29
$ws.Range("M136").Value = -5167.3638
$ws.Range("N136").Value = -228010.29
$ws.Range("H139").Value = 79767.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 79767.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 79767.75
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -90047.75
